$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 85.36364133333332
$ws.Range("H2").Value = 256.090924
$ws.Range("I2").Value = 0.832590152283795
$ws.Range("J2").Value = 0.8325901522837948
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 8.970048
$ws.Range("N2").Value = 26.910144
$ws.Range("O2").Value = 0.487108783009476
$ws.Range("P2").Value = 0.4871087830094759
$ws.Range("Q2").Value = 765.7159602147839
$ws.Range("R2").Value = 6891.443641933055
$ws.Range("S2").Value = 0.4055619758246336
$ws.Range("T2").Value = 0.4055619758246335

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 85.36364133333332
$ws.Range("H3").Value = 256.090924
$ws.Range("I3").Value = 0.832590152283795
$ws.Range("J3").Value = 0.8325901522837948
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 9.012070666666666
$ws.Range("N3").Value = 27.036212
$ws.Range("O3").Value = 0.489390778604016
$ws.Range("P3").Value = 0.489390778604016
$ws.Range("Q3").Value = 769.3031680599875
$ws.Range("R3").Value = 6923.728512539887
$ws.Range("S3").Value = 0.4074619428842026
$ws.Range("T3").Value = 0.4074619428842026

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 85.36364133333332
$ws.Range("H4").Value = 256.090924
$ws.Range("I4").Value = 0.832590152283795
$ws.Range("J4").Value = 0.8325901522837948
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.4327576666666667
$ws.Range("N4").Value = 1.298273
$ws.Range("O4").Value = 0.02350043838650813
$ws.Range("P4").Value = 0.02350043838650813
$ws.Range("Q4").Value = 36.94177024158355
$ws.Range("R4").Value = 332.475932174252
$ws.Range("S4").Value = 0.01956623357495874
$ws.Range("T4").Value = 0.01956623357495874

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 11.230072
$ws.Range("H5").Value = 33.690216
$ws.Range("I5").Value = 0.1095319647872954
$ws.Range("J5").Value = 0.1095319647872954
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 8.970048
$ws.Range("N5").Value = 26.910144
$ws.Range("O5").Value = 0.487108783009476
$ws.Range("P5").Value = 0.4871087830094759
$ws.Range("Q5").Value = 100.734284883456
$ws.Range("R5").Value = 906.608563951104
$ws.Range("S5").Value = 0.05335398206817624
$ws.Range("T5").Value = 0.05335398206817624

$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 11.230072
$ws.Range("H6").Value = 33.690216
$ws.Range("I6").Value = 0.1095319647872954
$ws.Range("J6").Value = 0.1095319647872954
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 9.012070666666666
$ws.Range("N6").Value = 27.036212
$ws.Range("O6").Value = 0.489390778604016
$ws.Range("P6").Value = 0.489390778604016
$ws.Range("Q6").Value = 101.2062024557547
$ws.Range("R6").Value = 910.855822101792
$ws.Range("S6").Value = 0.05360393352928217
$ws.Range("T6").Value = 0.05360393352928216

$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 11.230072
$ws.Range("H7").Value = 33.690216
$ws.Range("I7").Value = 0.1095319647872954
$ws.Range("J7").Value = 0.1095319647872954
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.4327576666666667
$ws.Range("N7").Value = 1.298273
$ws.Range("O7").Value = 0.02350043838650813
$ws.Range("P7").Value = 0.02350043838650813
$ws.Range("Q7").Value = 4.859899755218667
$ws.Range("R7").Value = 43.739097796968
$ws.Range("S7").Value = 0.002574049189837013
$ws.Range("T7").Value = 0.002574049189837013

$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 5.934092333333333
$ws.Range("H8").Value = 17.802277
$ws.Range("I8").Value = 0.05787788292890966
$ws.Range("J8").Value = 0.05787788292890966
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 8.970048
$ws.Range("N8").Value = 26.910144
$ws.Range("O8").Value = 0.487108783009476
$ws.Range("P8").Value = 0.4871087830094759
$ws.Range("Q8").Value = 53.229093066432
$ws.Range("R8").Value = 479.061837597888
$ws.Range("S8").Value = 0.02819282511666611
$ws.Range("T8").Value = 0.02819282511666611

$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 5.934092333333333
$ws.Range("H9").Value = 17.802277
$ws.Range("I9").Value = 0.05787788292890966
$ws.Range("J9").Value = 0.05787788292890966
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 9.012070666666666
$ws.Range("N9").Value = 27.036212
$ws.Range("O9").Value = 0.489390778604016
$ws.Range("P9").Value = 0.489390778604016
$ws.Range("Q9").Value = 53.47845945052489
$ws.Range("R9").Value = 481.306135054724
$ws.Range("S9").Value = 0.02832490219053118
$ws.Range("T9").Value = 0.02832490219053118

$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 5.934092333333333
$ws.Range("H10").Value = 17.802277
$ws.Range("I10").Value = 0.05787788292890966
$ws.Range("J10").Value = 0.05787788292890966
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.4327576666666667
$ws.Range("N10").Value = 1.298273
$ws.Range("O10").Value = 0.02350043838650813
$ws.Range("P10").Value = 0.02350043838650813
$ws.Range("Q10").Value = 2.568023951957889
$ws.Range("R10").Value = 23.112215567621
$ws.Range("S10").Value = 0.001360155621712372
$ws.Range("T10").Value = 0.001360155621712372
